$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1243.9073
$ws.Range("I137").Value = 986.2619
$ws.Range("J137").Value = 2145.6667
$ws.Range("K137").Value = 2958.7857
$ws.Range("L137").Value = 6437.000100000001
$ws.Range("M137").Value = -408.7856999999999
$ws.Range("N137").Value = -11537.0001

$ws.Range("H138").Value = 2012.4789
$ws.Range("I138").Value = 577.3461
$ws.Range("J138").Value = 5940.2104
$ws.Range("K138").Value = 1732.0383
$ws.Range("L138").Value = 17820.6312
$ws.Range("M138").Value = 3407.9617
$ws.Range("N138").Value = -28100.6312

$ws.Range("H141").Value = 1306.2037
$ws.Range("I141").Value = 866.0952
$ws.Range("J141").Value = 2846.5833
$ws.Range("K141").Value = 2598.2856
$ws.Range("L141").Value = 8539.749899999999
$ws.Range("M141").Value = 2581.7144
$ws.Range("N141").Value = -18899.7499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2337.5789
$ws.Range("I61").Value = 2333.9863
$ws.Range("J61").Value = 2425
$ws.Range("K61").Value = 2333.9863
$ws.Range("L61").Value = 2425
$ws.Range("M61").Value = -2121.9863
$ws.Range("N61").Value = -2849

$ws.Range("H74").Value = 8334423
$ws.Range("I74").Value = 922.6326
$ws.Range("K74").Value = 922.6326
$ws.Range("M74").Value = -48.63260000000002

$ws.Range("H77").Value = 8334423
$ws.Range("I77").Value = 922.6326
$ws.Range("K77").Value = 4613.163
$ws.Range("M77").Value = -245.1630000000005

$ws.Range("H97").Value = 1213.5883
$ws.Range("I97").Value = 1100.6923
$ws.Range("J97").Value = 1580.5
$ws.Range("K97").Value = 1100.6923
$ws.Range("L97").Value = 1580.5
$ws.Range("M97").Value = -604.6922999999999
$ws.Range("N97").Value = -2572.5

$ws.Range("H132").Value = 3944
$ws.Range("I132").Value = 2953.4211
$ws.Range("J132").Value = 5655
$ws.Range("K132").Value = 8860.263300000001
$ws.Range("L132").Value = 16965
$ws.Range("M132").Value = -6330.263300000001
$ws.Range("N132").Value = -22025

$ws.Range("H136").Value = 2337.5789
$ws.Range("I136").Value = 2333.9863
$ws.Range("J136").Value = 2425
$ws.Range("K136").Value = 7001.9589
$ws.Range("L136").Value = 7275
$ws.Range("M136").Value = -4451.9589
$ws.Range("N136").Value = -12375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 26123.6
$ws.Range("I20").Value = 1048
$ws.Range("K20").Value = 1048
$ws.Range("M20").Value = -801

$ws.Range("H86").Value = 14494243
$ws.Range("I86").Value = 19609252
$ws.Range("K86").Value = 19609252
$ws.Range("M86").Value = -19608129

$ws.Range("H89").Value = 14494243
$ws.Range("I89").Value = 19609252
$ws.Range("K89").Value = 98046260
$ws.Range("M89").Value = -98040644

$ws.Range("H134").Value = 4725.718
$ws.Range("I134").Value = 5862.423
$ws.Range("K134").Value = 17587.269
$ws.Range("M134").Value = -15052.269

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6496150
$ws.Range("I31").Value = 1727.6904
$ws.Range("J31").Value = 14289456
$ws.Range("K31").Value = 1727.6904
$ws.Range("L31").Value = 14289456
$ws.Range("M31").Value = -1432.6904
$ws.Range("N31").Value = -14290046

$ws.Range("H34").Value = 6496150
$ws.Range("I34").Value = 1727.6904
$ws.Range("J34").Value = 14289456
$ws.Range("K34").Value = 1727.6904
$ws.Range("L34").Value = 14289456
$ws.Range("M34").Value = -1525.6904
$ws.Range("N34").Value = -14289860

$ws.Range("H58").Value = 2565209.5
$ws.Range("I58").Value = 4065696.8
$ws.Range("J58").Value = 1877.25
$ws.Range("K58").Value = 4065696.8
$ws.Range("L58").Value = 1877.25
$ws.Range("M58").Value = -4065493.8
$ws.Range("N58").Value = -2283.25

$ws.Range("H105").Value = 16669048
$ws.Range("I105").Value = 19610204
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 19610204
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = -19608457
$ws.Range("N105").Value = -5994

$ws.Range("H122").Value = 2659.4167
$ws.Range("I122").Value = 2351.625
$ws.Range("J122").Value = 3275
$ws.Range("K122").Value = 7054.875
$ws.Range("L122").Value = 9825
$ws.Range("M122").Value = -4604.875
$ws.Range("N122").Value = -14725

$ws.Range("H132").Value = 2501301.2
$ws.Range("I132").Value = 3077884.8
$ws.Range("J132").Value = 2772.1333
$ws.Range("K132").Value = 9233654.399999999
$ws.Range("L132").Value = 8316.3999
$ws.Range("M132").Value = -9231124.399999999
$ws.Range("N132").Value = -13376.3999

$ws.Range("H134").Value = 8549351
$ws.Range("I134").Value = 11907458
$ws.Range("J134").Value = 1440.3636
$ws.Range("K134").Value = 35722374
$ws.Range("L134").Value = 4321.0908
$ws.Range("M134").Value = -35719839
$ws.Range("N134").Value = -9391.0908

$ws.Range("H136").Value = 2565209.5
$ws.Range("I136").Value = 4065696.8
$ws.Range("J136").Value = 1877.25
$ws.Range("K136").Value = 12197090.4
$ws.Range("L136").Value = 5631.75
$ws.Range("M136").Value = -12194540.4
$ws.Range("N136").Value = -10731.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 1744
$ws.Range("J80").Value = 1744
$ws.Range("L80").Value = 5232
$ws.Range("N80").Value = -7104

$ws.Range("H83").Value = 1744
$ws.Range("J83").Value = 1744
$ws.Range("L83").Value = 15696
$ws.Range("N83").Value = -25056

$ws.Range("H122").Value = 3847.9375
$ws.Range("I122").Value = 392.5
$ws.Range("J122").Value = 5921.2
$ws.Range("K122").Value = 3532.5
$ws.Range("L122").Value = 53290.8
$ws.Range("M122").Value = -1082.5
$ws.Range("N122").Value = -58190.8

$ws.Range("H131").Value = 1515886.8
$ws.Range("I131").Value = 2857601.2
$ws.Range("J131").Value = 1047.6129
$ws.Range("K131").Value = 8572803.600000001
$ws.Range("L131").Value = 3142.8387
$ws.Range("M131").Value = -8567763.600000001
$ws.Range("N131").Value = -13222.8387

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 9807461
$ws.Range("I132").Value = 15154631
$ws.Range("J132").Value = 4316.1665
$ws.Range("K132").Value = 45463893
$ws.Range("L132").Value = 12948.4995
$ws.Range("M132").Value = -45461363
$ws.Range("N132").Value = -18008.4995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1896.1786
$ws.Range("I7").Value = 1820.9474
$ws.Range("J7").Value = 2055
$ws.Range("K7").Value = 1820.9474
$ws.Range("L7").Value = 2055
$ws.Range("M7").Value = -1708.9474
$ws.Range("N7").Value = -2279

$ws.Range("H16").Value = 684.2857
$ws.Range("I16").Value = 465
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 465
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -295
$ws.Range("N16").Value = -2340

$ws.Range("H40").Value = 32261184
$ws.Range("I40").Value = 47621456
$ws.Range("J40").Value = 4614.9
$ws.Range("K40").Value = 47621456
$ws.Range("L40").Value = 4614.9
$ws.Range("M40").Value = -47621320
$ws.Range("N40").Value = -4886.9

$ws.Range("H82").Value = 369001.53
$ws.Range("I82").Value = 477597.44
$ws.Range("J82").Value = 115611.11
$ws.Range("K82").Value = 477597.44
$ws.Range("L82").Value = 115611.11
$ws.Range("M82").Value = -477236.44
$ws.Range("N82").Value = -116333.11

$ws.Range("H85").Value = 369001.53
$ws.Range("I85").Value = 477597.44
$ws.Range("J85").Value = 115611.11
$ws.Range("K85").Value = 477597.44
$ws.Range("L85").Value = 115611.11
$ws.Range("M85").Value = -476349.44
$ws.Range("N85").Value = -118107.11

$ws.Range("H126").Value = 1896.1786
$ws.Range("I126").Value = 1820.9474
$ws.Range("J126").Value = 2055
$ws.Range("K126").Value = 5462.8422
$ws.Range("L126").Value = 6165
$ws.Range("M126").Value = -2992.8422
$ws.Range("N126").Value = -11105

$ws.Range("H132").Value = 17634770
$ws.Range("I132").Value = 21832382
$ws.Range("J132").Value = 4799.6
$ws.Range("K132").Value = 65497146
$ws.Range("L132").Value = 14398.8
$ws.Range("M132").Value = -65494616
$ws.Range("N132").Value = -19458.8

$ws.Range("H136").Value = 5405.408
$ws.Range("I136").Value = 4494.1055
$ws.Range("J136").Value = 8553.546
$ws.Range("K136").Value = 13482.3165
$ws.Range("L136").Value = 25660.638
$ws.Range("M136").Value = -10932.3165
$ws.Range("N136").Value = -30760.638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 16209
$ws.Range("J86").Value = 16209
$ws.Range("L86").Value = 16209
$ws.Range("N86").Value = -18455

$ws.Range("H89").Value = 16209
$ws.Range("J89").Value = 16209
$ws.Range("L89").Value = 81045
$ws.Range("N89").Value = -92277

$ws.Range("H132").Value = 1112.2885
$ws.Range("I132").Value = 624.3143
$ws.Range("J132").Value = 2116.9412
$ws.Range("K132").Value = 1872.9429
$ws.Range("L132").Value = 6350.823600000001
$ws.Range("M132").Value = 657.0571
$ws.Range("N132").Value = -11410.8236

$ws.Range("H136").Value = 3971673.2
$ws.Range("I136").Value = 4663.409
$ws.Range("J136").Value = 8335384
$ws.Range("K136").Value = 13990.227
$ws.Range("L136").Value = 25006152
$ws.Range("M136").Value = -11440.227
$ws.Range("N136").Value = -25011252
